$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Artfynd")

# The observation that was in row 17 (Spillkråka / Dryocopus martius) and the
# observation that was in row 18 (Blåsippa / Hepatica nobilis) trade places.
# Row 17 also carries K/L/M/N (Aktivitet-related) cells that belong to the
# bird record, so those need to move down to row 18 along with it.

$scratchRow = 1000  # far away, unused row used as temporary holding area

# Columns whose cell is present (value or typed-empty) in both row 17 and
# row 18 and simply need to swap places.
$swapCols = @("A","B","D","E","F","G","H","Q","R")

foreach ($col in $swapCols) {
    $ws.Range($col + "17").Copy($ws.Range($col + $scratchRow))
}
foreach ($col in $swapCols) {
    $ws.Range($col + "18").Copy($ws.Range($col + "17"))
}
foreach ($col in $swapCols) {
    $ws.Range($col + $scratchRow).Copy($ws.Range($col + "18"))
    $ws.Range($col + $scratchRow).ClearContents()
}

# K, L, M, N only exist on row 17 before the edit (the bird sighting's
# activity info) and must end up on row 18 after the edit; row 17 loses them.
$moveCols = @("K", "L", "M", "N")
foreach ($col in $moveCols) {
    $ws.Range($col + "17").Copy($ws.Range($col + "18"))
    $ws.Range($col + "17").ClearContents()
}
